# Removed Test Case Inter-Dependency
#
# The ProductLoanInput sheet's "productname" / "shortname" values were
# hard-coded in a way that tied this test case's generated loan product
# name to another test case (...-VAR-INST). Give this sheet its own,
# independent product name / short name instead, and leave the selection
# and active-tab state the way the author left the workbook (parked on the
# ProductLoanOutput/verification sheet).

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# --- ProductLoanInput -------------------------------------------------
$wsInput.Activate()

# productname: no longer shares the "...-VAR-INST" name with the other test
$wsInput.Range("B1").Value = "2535-MS-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-1st"

# shortname: own literal short code instead of the shared numeric value
$wsInput.Range("B2").Value = "253e"

# Selection parked back at the top of the sheet
[void]$wsInput.Range("B1").Select()

# --- ProductLoanOutput --------------------------------------------------
# Cached verification value follows the renamed product
$wsOutput.Range("B1").Value = "2535-MS-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-1st"

# Leave the workbook with the output/verification sheet active & selected
$wsOutput.Activate()
[void]$wsOutput.Range("B1").Select()
